$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/replace header cells (column order: A,C,D stay; B,E,F,G change)
# New string values are introduced in the order: Resource, Timer, Student Number
$ws.Range("E1").Value = "Resource"
$ws.Range("F1").Value = "Timer"
$ws.Range("B1").Value = "Student Number"
$ws.Range("G1").Value = $null

# Update selection to B2
$ws.Range("B2").Select()
